$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'261.01"
$ws.Range("E2").Value = "'-0.06%"

$ws.Range("D3").Value = "'26.85"
$ws.Range("E3").Value = "'-1.64%"

$ws.Range("D4").Value = "'4.696"
$ws.Range("E4").Value = "'-0.30%"

$ws.Range("D5").Value = "'0.06215"
$ws.Range("E5").Value = "'2.33%"

$ws.Range("D6").Value = "'6.756"
$ws.Range("E6").Value = "'1.36%"

$ws.Range("E7").Value = "'0.47%"

$ws.Range("D8").Value = "'0.9125"
$ws.Range("E8").Value = "'-0.51%"

$ws.Range("D9").Value = "'0.1404"
$ws.Range("E9").Value = "'0.02%"

$ws.Range("D10").Value = "'0.04915"
$ws.Range("E10").Value = "'1.38%"

$ws.Range("D11").Value = "'0.07092"
$ws.Range("E11").Value = "'0.00%"

$ws.Range("D12").Value = "'0.03100"
$ws.Range("E12").Value = "'-1.29%"

$ws.Range("D13").Value = "'0.09046"

$ws.Range("D14").Value = "'0.001533"
$ws.Range("E14").Value = "'-0.15%"

$ws.Range("D15").Value = "'0.0006157"
$ws.Range("E15").Value = "'1.30%"

$ws.Range("D16").Value = "'0.006036"
$ws.Range("E16").Value = "'-0.04%"

$ws.Range("D17").Value = "'3.445"
$ws.Range("E17").Value = "'-0.12%"

$ws.Range("E18").Value = "'1.01%"

$ws.Range("D19").Value = "'2.145"
$ws.Range("E19").Value = "'-1.41%"

$ws.Range("E21").Value = "'1.07%"

$ws.Range("D22").Value = "'4.087"
$ws.Range("E22").Value = "'-0.30%"

$ws.Range("D23").Value = "'0.04232"
$ws.Range("E23").Value = "'-0.70%"

$ws.Range("D24").Value = "'0.001204"
$ws.Range("E24").Value = "'-1.27%"

$ws.Range("D25").Value = "'0.004077"
$ws.Range("E25").Value = "'4.23%"

$ws.Range("D40").Value = "'0.03955"
$ws.Range("E40").Value = "'2.10%"

$ws.Range("D41").Value = "'0.1112"
$ws.Range("E41").Value = "'-0.06%"

$ws.Range("D42").Value = "'0.004135"
$ws.Range("E42").Value = "'0.18%"

$ws.Range("E44").Value = "'-15.13%"

$ws.Range("D45").Value = "'0.00005162"
$ws.Range("E45").Value = "'-3.14%"

$ws.Range("D48").Value = "'0.2600"
$ws.Range("E48").Value = "'92.13%"
